# Import-export round trip: the app now serializes related child records
# (Destination / Travellers / Reviews / Attractions / Traveller) as raw JSON
# text instead of a human-readable comma-joined summary. Update the cached
# export values accordingly and widen the columns that now hold JSON text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Trips sheet
# ---------------------------------------------------------------------------
$trips = $wb.Worksheets.Item("Trips")

# Row 2 ("put" / probno putovanje) -> Destination/Travellers/Reviews as JSON
$trips.Range("F2").Value = '{"DestinationId":1,"City":"Paris","Country":"France"}'
$trips.Range("G2").Value = '[{"TravellerId":3,"Email":"mika@gmail.com","FirstName":"mika"}]'
$trips.Range("H2").Value = '[{"ReviewId":11,"Comment":"aq","Rating":10}]'

# Row 3 ("putovanje2") -> Destination JSON, plus previously-blank
# Travellers/Reviews cells now export as an explicit empty JSON array
$trips.Range("F3").Value = '{"DestinationId":5,"City":"nova Destinacija","Country":"novann"}'
$trips.Range("G3").Value = '[]'
$trips.Range("H3").Value = '[]'

# Row 4 ("w") -> same Paris destination JSON, new empty Travellers/Reviews
$trips.Range("F4").Value = '{"DestinationId":1,"City":"Paris","Country":"France"}'
$trips.Range("G4").Value = '[]'
$trips.Range("H4").Value = '[]'

# Columns grew a lot wider now that they hold JSON payloads
$trips.Columns.Item(6).ColumnWidth = 52.166666
$trips.Columns.Item(7).ColumnWidth = 53.666666
$trips.Columns.Item(8).ColumnWidth = 37.75

# ---------------------------------------------------------------------------
# Destinations sheet
# ---------------------------------------------------------------------------
$destinations = $wb.Worksheets.Item("Destinations")

# Row 2 (Paris) / Row 3 (New York) Attractions column -> JSON array
$destinations.Range("F2").Value = '[{"AttractionId":3,"Name":"Notre Dame Cathedral"},{"AttractionId":4,"Name":"Montmartre"}]'
$destinations.Range("F3").Value = '[{"AttractionId":2,"Name":"Central Park"},{"AttractionId":32,"Name":"Statue of Libertyy"}]'
# Row 4 (nova Destinacija) had no Attractions before; now exports "[]"
$destinations.Range("F4").Value = '[]'

$destinations.Columns.Item(6).ColumnWidth = 74.75

# ---------------------------------------------------------------------------
# Reviews sheet
# ---------------------------------------------------------------------------
$reviews = $wb.Worksheets.Item("Reviews")

# Row 2 Traveller column -> JSON instead of bare email
$reviews.Range("F2").Value = '{"TravellerId":3,"Email":"mika@gmail.com","FirstName":"mika"}'

$reviews.Columns.Item(6).ColumnWidth = 52.499999
